# Apply crypto price/volume updates scraped on Tue Apr 18 21:07:44 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the price cells whose new values look like plain numbers as Text first,
# so Excel keeps storing them as strings (matching the source data feed format)
# instead of auto-converting them into numeric cells.
$textForceCells = @("D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D17", "D18", "D19", "D20", "D21", "D22", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '30.391.18'
$ws.Range("E2").Value = '  +2.32%  '
$ws.Range("D3").Value = '2.094.05'
$ws.Range("E3").Value = '  -0.05%  '
$ws.Range("E4").Value = '  -0.97%  '
$ws.Range("D5").Value = '343.08'
$ws.Range("E5").Value = '  -0.52%  '
$ws.Range("E6").Value = '  -0.87%  '
$ws.Range("D7").Value = '0.5241'
$ws.Range("E7").Value = '  +1.58%  '
$ws.Range("D8").Value = '0.4422'
$ws.Range("E8").Value = '  +0.73%  '
$ws.Range("D9").Value = '54.49'
$ws.Range("E9").Value = '  +3.57%  '
$ws.Range("D10").Value = '0.09309'
$ws.Range("E10").Value = '  +0.67%  '
$ws.Range("D11").Value = '1.169'
$ws.Range("E11").Value = '  +0.00%  '
$ws.Range("D12").Value = '24.81'
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").Value = '8.607'
$ws.Range("E13").Value = '  +3.78%  '
$ws.Range("D14").Value = '6.898'
$ws.Range("E14").Value = '  +2.26%  '
$ws.Range("D15").Value = '2.043.95'
$ws.Range("E15").Value = '  -2.25%  '
$ws.Range("E16").Value = '  +1.99%  '
$ws.Range("D17").Value = '0.00001158'
$ws.Range("E17").Value = '  +0.72%  '
$ws.Range("D18").Value = '1.002'
$ws.Range("E18").Value = '  -0.85%  '
$ws.Range("D19").Value = '21.13'
$ws.Range("E19").Value = '  +1.45%  '
$ws.Range("D20").Value = '0.06652'
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("D21").Value = '6.325'
$ws.Range("E21").Value = '  +2.23%  '
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.83%  '
$ws.Range("D23").Value = '30.394.56'
$ws.Range("E23").Value = '  +2.23%  '
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").Value = '2.302'
$ws.Range("E25").Value = '  -0.72%  '
$ws.Range("D26").Value = '21.82'
$ws.Range("E26").Value = '  -0.41%  '
$ws.Range("D27").Value = '162.84'
$ws.Range("E27").Value = '  +0.73%  '
$ws.Range("D28").Value = '2.506'
$ws.Range("E28").Value = '  -0.38%  '
$ws.Range("D29").Value = '133.28'
$ws.Range("E29").Value = '  +0.22%  '
$ws.Range("D30").Value = '1.140'
$ws.Range("E30").Value = '  +0.80%  '
$ws.Range("D31").Value = '1.682'
$ws.Range("E31").Value = '  +1.95%  '
$ws.Range("D32").Value = '0.1046'
$ws.Range("E32").Value = '  -0.50%  '
$ws.Range("D33").Value = '6.822'
$ws.Range("E33").Value = '  +9.74%  '
$ws.Range("D34").Value = '6.254'
$ws.Range("E34").Value = '  +1.42%  '
$ws.Range("D35").Value = '3.866'
$ws.Range("E35").Value = '  -1.74%  '
$ws.Range("D36").Value = '10.16'
$ws.Range("E36").Value = '  -0.97%  '
$ws.Range("D37").Value = '0.02637'
$ws.Range("E37").Value = '  +2.72%  '
$ws.Range("D38").Value = '0.06841'
$ws.Range("E38").Value = '  +2.10%  '
$ws.Range("D39").Value = '0.6991'
$ws.Range("E39").Value = '  +1.65%  '
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").Value = '12.56'
$ws.Range("E40").Value = '  +0.95%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '1.343'
$ws.Range("E41").Value = '  +2.78%  '
$ws.Range("D42").Value = '0.2211'
$ws.Range("E42").Value = '  -0.61%  '
$ws.Range("D43").Value = '0.6816'
$ws.Range("E43").Value = '  +2.25%  '
$ws.Range("D44").Value = '14.36'
$ws.Range("E44").Value = '  +0.62%  '
$ws.Range("E45").Value = '  +1.31%  '
$ws.Range("D46").Value = '0.9991'
$ws.Range("E46").Value = '  -0.86%  '
$ws.Range("D47").Value = '1.370'
$ws.Range("E47").Value = '  +18.05%  '
$ws.Range("D48").Value = '3.632'
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("D49").Value = '0.00000000352'
$ws.Range("E49").Value = '  -0.57%  '
$ws.Range("D50").Value = '1.213'
$ws.Range("E50").Value = '  +8.45%  '
$ws.Range("E51").Value = '  -0.19%  '
